# Update cryptocurrency price/volume data per latest GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D hold price text that can look like plain decimals (e.g. "42.36").
# Force those specific cells to Text format first so Excel keeps the exact
# string instead of silently converting it to a floating point number.
$textCells = @('D5', 'D8', 'D10', 'D11', 'D14', 'D18', 'D20', 'D22', 'D25', 'D26', 'D27', 'D29', 'D33', 'D34', 'D37', 'D38', 'D39', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '35.101.76'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').Value = '1.857.19'
$ws.Range('E3').Value = '  +3.10%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '237.08'
$ws.Range('E5').Value = '  +3.50%  '
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('D8').Value = '42.36'
$ws.Range('E8').Value = '  +7.55%  '
$ws.Range('D10').Value = '0.0696'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('D11').Value = '0.0991'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '2.126.01'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.873.85'
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '11.44'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('E15').Value = '  +3.05%  '
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').Value = '35.059.37'
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').Value = '70.34'
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('D20').Value = '240.93'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('E21').Value = '  +3.07%  '
$ws.Range('D22').Value = '4.76'
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('D25').Value = '171.25'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').Value = '1.89'
$ws.Range('E26').Value = '  +26.89%  '
$ws.Range('D27').Value = '7.95'
$ws.Range('E27').Value = '  +2.63%  '
$ws.Range('E28').Value = '  +3.04%  '
$ws.Range('D29').Value = '0.125'
$ws.Range('E29').Value = '  +2.60%  '
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '4.02'
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('D34').Value = '2.02'
$ws.Range('E34').Value = '  +13.02%  '
$ws.Range('E35').Value = '  +23.06%  '
$ws.Range('E36').Value = '  +5.77%  '
$ws.Range('D37').Value = '0.784'
$ws.Range('E37').Value = '  +13.15%  '
$ws.Range('D38').Value = '1.09'
$ws.Range('E38').Value = '  +13.46%  '
$ws.Range('D39').Value = '91.73'
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('E40').Value = '  +7.26%  '
$ws.Range('D41').Value = '1.353.82'
$ws.Range('E41').Value = '  +2.32%  '
$ws.Range('D42').Value = '14.87'
$ws.Range('E42').Value = '  +4.87%  '
$ws.Range('D43').Value = '2.35'
$ws.Range('E43').Value = '  +6.08%  '
$ws.Range('D44').Value = '12.79'
$ws.Range('E44').Value = '  +57.12%  '
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('D46').Value = '2.75'
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('D47').Value = '0.0551'
$ws.Range('E47').Value = '  +7.60%  '
$ws.Range('D48').Value = '6.46'
$ws.Range('E48').Value = '  +5.29%  '
$ws.Range('D49').Value = '2.038.91'
$ws.Range('E49').Value = '  +2.77%  '
$ws.Range('E50').Value = '  +3.42%  '
$ws.Range('E51').Value = '  +18.20%  '
